# Add a new booking row (row 7) to the Bookings sheet, mirroring the
# existing "walk-in" style rows (e.g. row 6) where the date-like fields
# are stored as literal text rather than being auto-converted to Excel
# date serial numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bookings")

# --- Plain text / identifier columns -------------------------------------
$ws.Range("A7").Value = "SNOW-401902"

# --- Date-like text columns -------------------------------------------
# Force text storage (NumberFormat "@") so Excel doesn't reinterpret the
# string as a date serial, then restore the default "Normal" style so the
# cell doesn't carry a lingering custom number-format style index.
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2/25/2026"
$ws.Range("B7").Style = "Normal"

$ws.Range("C7").Value = "local"
$ws.Range("D7").Value = "l@l.com"
$ws.Range("E7").Value = "local1231"

# --- Numeric columns -----------------------------------------------------
$ws.Range("F7").Value = 1

$ws.Range("G7").Value = "Family Ski Package"

$ws.Range("H7").Value = 32000
$ws.Range("I7").Value = 32000

$ws.Range("J7").Value = "Confirmed"

# --- Second date-like text column ----------------------------------------
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "2/23/2026"
$ws.Range("K7").Style = "Normal"

# --- Special Requests: empty string (not a blank/absent cell) ------------
# A lone apostrophe is Excel's quote-prefix for an empty text entry; this
# produces an actual empty-string cell (matching the existing L6 cell)
# rather than simply clearing/removing the cell.
$ws.Range("L7").Value = "'"
$ws.Range("L7").Style = "Normal"

Write-Output "Row 7 booking added."
